# Update Okanogan-related habitat quality scores on rows 26-29.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 26: J26 3->5, K26 5->(blank), N26 (blank)->5
$ws.Range("J26").Value = 5
$ws.Range("K26").ClearContents()
$ws.Range("N26").Value = 5

# Row 27: J27 3->5, K27 5->3, N27 (blank)->5
$ws.Range("J27").Value = 5
$ws.Range("K27").Value = 3
$ws.Range("N27").Value = 5

# Row 28: J28 3->5, N28 (blank)->5 (K28 unchanged)
$ws.Range("J28").Value = 5
$ws.Range("N28").Value = 5

# Row 29: J29 3->5, K29 3->(blank), N29 (blank)->3
$ws.Range("J29").Value = 5
$ws.Range("K29").ClearContents()
$ws.Range("N29").Value = 3
